$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "35.364.35"
$ws.Range("E2").Value = "  +1.62%  "

$ws.Range("D3").Value = "1.895.61"
$ws.Range("E3").Value = "  +1.73%  "

$ws.Range("E4").Value = "  +0.01%  "

$ws.Range("D5").Value = "'246.78"
$ws.Range("E5").Value = "  +1.11%  "

$ws.Range("D6").Value = "'0.691"
$ws.Range("E6").Value = "  +3.00%  "

$ws.Range("E7").Value = "  +0.02%  "

$ws.Range("D8").Value = "'42.96"
$ws.Range("E8").Value = "  +2.28%  "

$ws.Range("E9").Value = "  +5.77%  "

$ws.Range("D10").Value = "'55.91"
$ws.Range("E10").Value = "  +10.15%  "

$ws.Range("E11").Value = "  +2.65%  "

$ws.Range("E12").Value = "  +1.56%  "

$ws.Range("D13").Value = "'13.98"
$ws.Range("E13").Value = "  +9.27%  "

$ws.Range("E14").Value = "  +12.55%  "

$ws.Range("D15").Value = "2.171.47"
$ws.Range("E15").Value = "  +1.96%  "

$ws.Range("D16").Value = "'5.01"
$ws.Range("E16").Value = "  +4.68%  "

$ws.Range("D17").Value = "1.873.64"
$ws.Range("E17").Value = "  +0.53%  "

$ws.Range("D18").Value = "35.435.59"
$ws.Range("E18").Value = "  +1.92%  "

$ws.Range("D19").Value = "'73.54"
$ws.Range("E19").Value = "  +2.48%  "

$ws.Range("D20").Value = "0.0₃0827"
$ws.Range("E20").Value = "  +2.87%  "

$ws.Range("D21").Value = "'244.61"
$ws.Range("E21").Value = "  +1.22%  "

$ws.Range("D22").Value = "'12.92"
$ws.Range("E22").Value = "  +3.50%  "

$ws.Range("E23").Value = "  +7.98%  "

$ws.Range("D24").Value = "'2.67"
$ws.Range("E24").Value = "  +7.62%  "

$ws.Range("D26").Value = "'2.18"
$ws.Range("E26").Value = "  +1.34%  "

$ws.Range("D27").Value = "'166.50"
$ws.Range("E27").Value = "  +2.37%  "

$ws.Range("E28").Value = "  +3.10%  "

$ws.Range("E29").Value = "  +2.12%  "

$ws.Range("D30").Value = "'0.128"
$ws.Range("E30").Value = "  +2.24%  "

$ws.Range("D31").Value = "'0.0605"
$ws.Range("E31").Value = "  +6.62%  "

$ws.Range("D32").Value = "'4.35"
$ws.Range("E32").Value = "  +5.24%  "

$ws.Range("D33").Value = "'4.23"
$ws.Range("E33").Value = "  +3.55%  "

$ws.Range("D34").Value = "'1.87"
$ws.Range("E34").Value = "  +24.40%  "

$ws.Range("E35").Value = "  -0.04%  "

$ws.Range("E36").Value = "  -13.95%  "

$ws.Range("D37").Value = "'0.852"
$ws.Range("E37").Value = "  +3.83%  "

$ws.Range("E38").Value = "  +2.07%  "

$ws.Range("D39").Value = "'0.0714"
$ws.Range("E39").Value = "  +8.19%  "

$ws.Range("D40").Value = "'0.0224"
$ws.Range("E40").Value = "  +7.30%  "

$ws.Range("D41").Value = "'99.05"
$ws.Range("E41").Value = "  +1.99%  "

$ws.Range("D42").Value = "'17.01"
$ws.Range("E42").Value = "  +0.79%  "

$ws.Range("E43").Value = "  +1.97%  "

$ws.Range("D44").Value = "1.338.89"
$ws.Range("E44").Value = "  +5.03%  "

$ws.Range("D45").Value = "'13.54"
$ws.Range("E45").Value = "  +14.25%  "

$ws.Range("E46").Value = "  +3.99%  "

$ws.Range("D47").Value = "'0.0809"
$ws.Range("E47").Value = "  -3.78%  "

$ws.Range("E48").Value = "  +0.98%  "

$ws.Range("E49").Value = "  +0.35%  "

$ws.Range("E50").Value = "  +2.38%  "

$ws.Range("B51").Value = "MultiversX"
$ws.Range("C51").Value = "https://coinranking.com/coin/omwkOTglq+multiversx-egld"
$ws.Range("D51").Value = "'42.57"
$ws.Range("E51").Value = "  +1.14%  "
